$wb = $excel.ActiveWorkbook

# --- About sheet: update the "last updated" date in C1 (2024-03-18 -> 2024-04-05) ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45387

# --- BAU Emissions sheet: bulk label rename " : NoSettings" -> " : test" for rows 4-278 (col A) ---
$bau = $wb.Worksheets.Item("BAU Emissions")
for ($r = 4; $r -le 278; $r++) {
    $cell = $bau.Cells.Item($r, 1)
    $old = $cell.Value2
    $new = $old -replace " : NoSettings$", " : test"
    $cell.Value = $new
}

# --- BAU Emissions sheet: update forecast values on row 94 (columns M:AE) ---
$bau.Range("M94").Value = 1001080
$bau.Range("N94").Value = 2002150
$bau.Range("O94").Value = 3003230
$bau.Range("P94").Value = 4004300
$bau.Range("Q94").Value = 5005380
$bau.Range("R94").Value = 5005380
$bau.Range("S94").Value = 5005380
$bau.Range("T94").Value = 5005380
$bau.Range("U94").Value = 5005380
$bau.Range("V94").Value = 5005380
$bau.Range("W94").Value = 5005380
$bau.Range("X94").Value = 5005380
$bau.Range("Y94").Value = 5005380
$bau.Range("Z94").Value = 5005380
$bau.Range("AA94").Value = 5005380
$bau.Range("AB94").Value = 5005380
$bau.Range("AC94").Value = 5005380
$bau.Range("AD94").Value = 5005380
$bau.Range("AE94").Value = 5005380

# --- Move the active tab from "Current and Planned Capacity" to "About" ---
# and restore the BAU Emissions selection/scroll state before we leave it.
$bau.Range("A30:AE280").Select()

$about.Activate()
$about.Range("E29").Select()
